# new Madigan bike hours
# Update Riders (C) and Average (D) columns for rows 2-8 on the Ridership sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

$ws.Range("C2").Value = 201
$ws.Range("D2").Value = 226.04

$ws.Range("C3").Value = 195
$ws.Range("D3").Value = 218.83

$ws.Range("C4").Value = 262
$ws.Range("D4").Value = 216.27

$ws.Range("C5").Value = 234
$ws.Range("D5").Value = 234.86

$ws.Range("C6").Value = 228
$ws.Range("D6").Value = 238.6

$ws.Range("C7").Value = 150
$ws.Range("D7").Value = 114.06

$ws.Range("C8").Value = 107
$ws.Range("D8").Value = 91.06999999999999
